$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 443.33334
$ws.Range("I12").Value = 1030
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 1030
$ws.Range("L12").Value = 150
$ws.Range("M12").Value = -860
$ws.Range("N12").Value = -490
$ws.Range("H33").Value = 125.2
$ws.Range("I33").Value = 136.875
$ws.Range("K33").Value = 136.875
$ws.Range("M33").Value = 92.125
$ws.Range("H51").Value = 6333.1665
$ws.Range("J51").Value = 6599.8
$ws.Range("L51").Value = 6599.8
$ws.Range("N51").Value = -7567.8
$ws.Range("H129").Value = 915.7857
$ws.Range("J129").Value = 887.28204
$ws.Range("L129").Value = 2661.84612
$ws.Range("N129").Value = -12661.84612
$ws.Range("H132").Value = 1344.9048
$ws.Range("I132").Value = 1265.6875
$ws.Range("K132").Value = 3797.0625
$ws.Range("M132").Value = -1267.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1292370
$ws.Range("I2").Value = 1789245.9
$ws.Range("J2").Value = 492.4
$ws.Range("K2").Value = 1789245.9
$ws.Range("L2").Value = 492.4
$ws.Range("M2").Value = -1789132.9
$ws.Range("N2").Value = -718.4
$ws.Range("H32").Value = 3086.5574
$ws.Range("I32").Value = 2063.1091
$ws.Range("K32").Value = 2063.1091
$ws.Range("M32").Value = -1776.1091
$ws.Range("H61").Value = 3846.25
$ws.Range("I61").Value = 2676.25
$ws.Range("K61").Value = 2676.25
$ws.Range("M61").Value = -2464.25
$ws.Range("H110").Value = 2649.818
$ws.Range("I110").Value = 1517.4286
$ws.Range("J110").Value = 4631.5
$ws.Range("K110").Value = 1517.4286
$ws.Range("L110").Value = 4631.5
$ws.Range("M110").Value = 527.5714
$ws.Range("N110").Value = -8721.5
$ws.Range("H116").Value = 1292370
$ws.Range("I116").Value = 1789245.9
$ws.Range("J116").Value = 492.4
$ws.Range("K116").Value = 1789245.9
$ws.Range("L116").Value = 492.4
$ws.Range("M116").Value = -1786951.9
$ws.Range("N116").Value = -5080.4
$ws.Range("H132").Value = 1683.9375
$ws.Range("I132").Value = 1084.5714
$ws.Range("K132").Value = 3253.7142
$ws.Range("M132").Value = -723.7142000000003
$ws.Range("H135").Value = 19214.5
$ws.Range("J135").Value = 19214.5
$ws.Range("L135").Value = 19214.5
$ws.Range("N135").Value = -29354.5
$ws.Range("H136").Value = 3846.25
$ws.Range("I136").Value = 2676.25
$ws.Range("K136").Value = 8028.75
$ws.Range("M136").Value = -5478.75
$ws.Range("H139").Value = 46905
$ws.Range("J139").Value = 46905
$ws.Range("L139").Value = 46905
$ws.Range("N139").Value = -57185

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1292370
$ws.Range("I3").Value = 1789245.9
$ws.Range("J3").Value = 492.4
$ws.Range("K3").Value = 1789245.9
$ws.Range("L3").Value = 492.4
$ws.Range("M3").Value = -1789131.9
$ws.Range("N3").Value = -720.4
$ws.Range("H108").Value = 94981.5
$ws.Range("J108").Value = 94981.5
$ws.Range("L108").Value = 94981.5
$ws.Range("N108").Value = -102661.5
$ws.Range("H134").Value = 10371.069
$ws.Range("I134").Value = 10293.458
$ws.Range("K134").Value = 30880.374
$ws.Range("M134").Value = -28345.374

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 770.3333
$ws.Range("J16").Value = 999.5
$ws.Range("L16").Value = 999.5
$ws.Range("N16").Value = -1573.5
$ws.Range("H58").Value = 1611875.9
$ws.Range("I58").Value = 3624313
$ws.Range("J58").Value = 1926.0667
$ws.Range("K58").Value = 3624313
$ws.Range("L58").Value = 1926.0667
$ws.Range("M58").Value = -3624110
$ws.Range("N58").Value = -2332.0667
$ws.Range("H113").Value = 770.3333
$ws.Range("J113").Value = 999.5
$ws.Range("L113").Value = 999.5
$ws.Range("N113").Value = -5339.5
$ws.Range("H132").Value = 2344.2727
$ws.Range("I132").Value = 1179.3636
$ws.Range("J132").Value = 3509.182
$ws.Range("K132").Value = 3538.0908
$ws.Range("L132").Value = 10527.546
$ws.Range("M132").Value = -1008.0908
$ws.Range("N132").Value = -15587.546
$ws.Range("H134").Value = 999.7273
$ws.Range("I134").Value = 999.7
$ws.Range("J134").Value = 1000
$ws.Range("K134").Value = 2999.1
$ws.Range("L134").Value = 3000
$ws.Range("M134").Value = -464.1000000000004
$ws.Range("N134").Value = -8070
$ws.Range("H136").Value = 1611875.9
$ws.Range("I136").Value = 3624313
$ws.Range("J136").Value = 1926.0667
$ws.Range("K136").Value = 10872939
$ws.Range("L136").Value = 5778.2001
$ws.Range("M136").Value = -10870389
$ws.Range("N136").Value = -10878.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 259.16666
$ws.Range("J92").Value = 271
$ws.Range("L92").Value = 813
$ws.Range("N92").Value = -3309
$ws.Range("H137").Value = 3598.0625
$ws.Range("I137").Value = 2835
$ws.Range("J137").Value = 3852.4167
$ws.Range("K137").Value = 8505
$ws.Range("L137").Value = 11557.2501
$ws.Range("M137").Value = -3405
$ws.Range("N137").Value = -21757.2501

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2560.375
$ws.Range("I102").Value = 2532.2856
$ws.Range("K102").Value = 2532.2856
$ws.Range("M102").Value = -910.2856000000002
$ws.Range("H109").Value = 17523.334
$ws.Range("J109").Value = 17523.334
$ws.Range("L109").Value = 17523.334
$ws.Range("N109").Value = -19603.334
$ws.Range("H113").Value = 1036
$ws.Range("I113").Value = 804
$ws.Range("K113").Value = 804
$ws.Range("M113").Value = 1366
$ws.Range("H122").Value = 1368.6
$ws.Range("I122").Value = 1020.3077
$ws.Range("K122").Value = 3060.9231
$ws.Range("M122").Value = -610.9231
$ws.Range("H132").Value = 1041933.1
$ws.Range("I132").Value = 1242586.6
$ws.Range("K132").Value = 3727759.8
$ws.Range("M132").Value = -3725229.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2854.8333
$ws.Range("J16").Value = 1966.2222
$ws.Range("L16").Value = 1966.2222
$ws.Range("N16").Value = -2306.2222
$ws.Range("H40").Value = 4588.5557
$ws.Range("J40").Value = 9685.286
$ws.Range("L40").Value = 9685.286
$ws.Range("N40").Value = -9957.286
$ws.Range("H61").Value = 2831
$ws.Range("I61").Value = 2618.4546
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2618.4546
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2416.4546
$ws.Range("N61").Value = -4404
$ws.Range("H93").Value = 485.83334
$ws.Range("I93").Value = 403.33334
$ws.Range("K93").Value = 403.33334
$ws.Range("M93").Value = 844.66666
$ws.Range("H113").Value = 2831
$ws.Range("I113").Value = 2618.4546
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2618.4546
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -448.4546
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 8756.556
$ws.Range("I122").Value = 6467.3335
$ws.Range("K122").Value = 19402.0005
$ws.Range("M122").Value = -16952.0005
$ws.Range("H133").Value = 89000
$ws.Range("J133").Value = 89000
$ws.Range("L133").Value = 89000
$ws.Range("N133").Value = -94060
$ws.Range("H136").Value = 3514.2942
$ws.Range("I136").Value = 1678.1428
$ws.Range("J136").Value = 4799.6
$ws.Range("K136").Value = 5034.428400000001
$ws.Range("L136").Value = 14398.8
$ws.Range("M136").Value = -2484.428400000001
$ws.Range("N136").Value = -19498.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1856.1666
$ws.Range("I81").Value = 1752.6666
$ws.Range("K81").Value = 3505.3332
$ws.Range("M81").Value = -2444.3332
$ws.Range("H84").Value = 1856.1666
$ws.Range("I84").Value = 1752.6666
$ws.Range("K84").Value = 17526.666
$ws.Range("M84").Value = -12222.666
$ws.Range("H132").Value = 2401.7778
$ws.Range("I132").Value = 1968.8235
$ws.Range("J132").Value = 3137.8
$ws.Range("K132").Value = 5906.470499999999
$ws.Range("L132").Value = 9413.400000000001
$ws.Range("M132").Value = -3376.470499999999
$ws.Range("N132").Value = -14473.4
